{"js": "// Apply the \"major changes and implementations\" edits to the interim report.\n//\n// Semantic changes made by this script:\n//   1. Team members line: complete the two truncated student IDs\n//      \"M13   \" -> \"M1341\" (Vasilis Skourtis) and \"M13    \" -> \"M1364\"\n//      (Evaggelos Karageorgos), and collapse the surrounding runs.\n//   2. \"...fallen behind with our scheduling.\" -> \"...fallen behind with\n//      our schedule.\" (dropped \"ing\").\n//   3. \"...that is why you are asking for a deadline extension\" ->\n//      \"...that is why we are asking for a deadline extension\".\n//   4. The stray \"_GoBack\" bookmark (Word's \"last edit position\" marker)\n//      is relocated from the \"skype meeting every other week\" sentence to\n//      the GitHub repository URL later in the document, reflecting where\n//      the author's final edit actually happened.\n//\n// The four hyperlink URLs (52.17.140.15:8080, the ELB 8080/8443 endpoints\n// and the GitHub link) keep their exact visible text - in the source\n// revision they only end up split into additional runs with no textual\n// change, so there is nothing to edit there from the API's point of view.\n\nconst doc = context.document;\nconst body = doc.body;\n\nasync function replaceOnce(searchText, replaceText) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(replaceText, \"Replace\");\n  await context.sync();\n}\n\n// 1. Complete the two student IDs in the \"Team members:\" line.\nawait replaceOnce(\"SKOURTIS (M13   )\", \"SKOURTIS (M1341)\");\nawait replaceOnce(\"KARAGEORGOS (M13    )\", \"KARAGEORGOS (M1364)\");\n\n// 2. \"scheduling\" -> \"schedule\".\nawait replaceOnce(\n  \"fallen behind with our scheduling\",\n  \"fallen behind with our schedule\"\n);\n\n// 3. \"you are asking\" -> \"we are asking\".\nawait replaceOnce(\n  \"that is why you are asking for a deadline extension\",\n  \"that is why we are asking for a deadline extension\"\n);\n\n// 4. Move the \"_GoBack\" bookmark from the \"every other week\" sentence to\n//    the GitHub URL near the end of the document.\nconst oldBookmark = doc.getBookmarkRangeOrNullObject(\"_GoBack\");\nawait context.sync();\nif (!oldBookmark.isNullObject) {\n  doc.deleteBookmark(\"_GoBack\");\n}\n\n// Target the same spot the diff shows (right after \"std080\", inside the\n// GitHub link); the engine snaps the collapsed bookmark range to the\n// nearest run boundary, which here is the end of the hyperlink run.\nconst githubResults = body.search(\"std080\", { matchCase: true });\ngithubResults.load(\"text\");\nawait context.sync();\nif (githubResults.items.length > 0) {\n  const githubRange = githubResults.items[0];\n  const collapsed = githubRange.getRange(\"End\");\n  collapsed.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Apply the \"major changes and implementations\" edits to the interim report.\n#\n# Semantic changes made by this script:\n#   1. Team members line: complete the two truncated student IDs\n#      \"M13   \" -> \"M1341\" (Vasilis Skourtis) and \"M13    \" -> \"M1364\"\n#      (Evaggelos Karageorgos).\n#   2. \"...fallen behind with our scheduling.\" -> \"...fallen behind with\n#      our schedule.\" (dropped \"ing\").\n#   3. \"...that is why you are asking for a deadline extension\" ->\n#      \"...that is why we are asking for a deadline extension\".\n#   4. The stray \"_GoBack\" bookmark (Word's \"last edit position\" marker)\n#      is relocated from the \"skype meeting every other week\" sentence to\n#      the GitHub repository URL later in the document, reflecting where\n#      the author's final edit actually happened.\n#\n# The four hyperlink URLs (52.17.140.15:8080, the ELB 8080/8443 endpoints\n# and the GitHub link) keep their exact visible text - nothing to change\n# there from a content point of view.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($searchText, $replaceText) {\n    $range = $d.Content\n    $found = $range.Find.Execute(\n        $searchText,  # FindText\n        $true,        # MatchCase\n        $false,       # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        1,            # Wrap (wdFindContinue)\n        $false,       # Format\n        $replaceText, # ReplaceWith\n        2             # Replace (wdReplaceOne)\n    )\n    if (-not $found) {\n        throw \"Text not found: $searchText\"\n    }\n}\n\n# 1. Complete the two student IDs in the \"Team members:\" line.\nReplace-Text \"SKOURTIS (M13   )\" \"SKOURTIS (M1341)\"\nReplace-Text \"KARAGEORGOS (M13    )\" \"KARAGEORGOS (M1364)\"\n\n# 2. \"scheduling\" -> \"schedule\".\nReplace-Text \"fallen behind with our scheduling\" \"fallen behind with our schedule\"\n\n# 3. \"you are asking\" -> \"we are asking\".\nReplace-Text \"that is why you are asking for a deadline extension\" \"that is why we are asking for a deadline extension\"\n\n# 4. Move the \"_GoBack\" bookmark from the \"every other week\" sentence to\n#    the GitHub URL near the end of the document.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Target the same spot the diff shows (right after \"std080\", inside the\n# GitHub link); the engine snaps the collapsed bookmark range to the\n# nearest run boundary, which here is the end of the hyperlink run.\n$ghRange = $d.Content\n$ghFound = $ghRange.Find.Execute(\"std080\")\nif ($ghFound) {\n    $ghRange.Collapse(0)  # wdCollapseEnd\n    $d.Bookmarks.Add(\"_GoBack\", $ghRange)\n}\n"}
